$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Pasteles Tere'
$ws.Range("B2").Value = 'https://fontsinuse.com/uses/37611/pasteles-tere?utm_source=feedburner&utm_medium=feed&utm_campaign=Feed%3A+FontsInUseAll+%28Fonts+In+Use%29'

$ws.Range("A3").Value = 'Bu Besinler Akciğerleri Güçlendiriyor!'
$ws.Range("B3").Value = 'https://www.hurriyet.com.tr/mahmure/galeri-bu-besinler-akcigerleri-guclendiriyor-41695361'

$ws.Range("A4").Value = 'Acusa panista de corrupción... a panista'
$ws.Range("B4").Value = 'https://www.reforma.com/libre/acceso/accesofb.htm?urlredirect=/acusa-panista-de-corrupcion-a-panista/ar2084120'

$ws.Range("A5").Value = 'Qué es el último teorema de Fermat y por qué los matemáticos demoraron 3 siglos en resolverlo'
$ws.Range("B5").Value = 'https://www.bbc.com/mundo/noticias-55412805'

$ws.Range("A6").Value = '¿Cómo sería México si Estados Unidos no se hubiera apropiado de más de la mitad de su territorio en el siglo XIX?'
$ws.Range("B6").Value = 'https://www.bbc.com/mundo/noticias-55151922'

$ws.Range("A7").Value = 'Cómo se rige la Antártida, quién reclama su soberanía y por qué despierta tanto interés'
$ws.Range("B7").Value = 'https://www.bbc.com/mundo/noticias-internacional-55108222'

$ws.Range("A8").Value = 'Mantecadas al tren, las vendedoras que hicieron viajera la dulce tradición de Astorga (León)'
$ws.Range("B8").Value = 'https://www.eldiario.es/castilla-y-leon/provincias/leon/mantecadas-tren-vendedoras-hicieron-viajera-dulce-tradicion-astorga-leon_1_6623637.html'

$ws.Range("A9").Value = 'Tere otunun faydaları nelerdir? İç organları arındıran tere otu nasıl tüketilir?'
$ws.Range("B9").Value = 'https://www.haber7.com/saglik/haber/2992834-tere-otunun-faydalari-nelerdir-ic-organlari-arindiran-tere-otu-nasil-tuketilir'

$ws.Range("A10").Value = 'Have bonded well with Akash Mukherjee: Aishwarya Raj Bhakuni'
$ws.Range("B10").Value = 'https://www.santabanta.com/bollywood/150620/have-bonded-well-with-akash-mukherjee-aishwarya-raj-bhakuni/'

$ws.Range("A11").Value = 'Topul dezvoltatorilor de software din România, domeniu în care lucrează peste 130.000 de IT-işt'
$ws.Range("B11").Value = 'https://economie.hotnews.ro/stiri-it-24454572-topul-dezvoltatorilor-software-din-romania-domeniu-care-lucreaza-peste-130-000-ist.htm'

$ws.Range("A12").Value = 'Noi finanțări de 1-3 milioane Euro pentru startup-uri IT românști'
$ws.Range("B12").Value = 'https://economie.hotnews.ro/stiri-eurofonduri-24454593-noi-finantari-1-3-milioane-euro-pentru-startup-uri-romansti.htm'

$ws.Range("A13").Value = 'O familie din Iasi s-a mutat în Zanzibar din cauza restricţiilor Covid. "E mai ieftin aici de trăit"'
$ws.Range("B13").Value = 'https://www.hotnews.ro/stiri-coronavirus-24453771-familie-din-iasi-mutat-zanzibar-din-cauza-restrictiilor-covid-mai-ieftin-aici-trait.htm'

$ws.Range("A14").Value = '​VIDEO Real Madrid, învinsă pentru a doua oară de Șahtior Donețk (2-0) / RB Salzburg, victorie în Rusia (3-1 vs Lokomotiv Moscova)'
$ws.Range("B14").Value = 'https://sport.hotnews.ro/stiri-fotbal-24453467-video-sahtior-donetk-invinge-pentru-doua-oara-real-madrid-2-0-salzburg-victorie-rusia-3-1-lokomotiv-moscova.htm'

$ws.Range("A15").Value = 'Mircea Rednic este noul antrenor al echipei FC Viitorul'
$ws.Range("B15").Value = 'https://sport.hotnews.ro/stiri-fotbal-24453452-mircea-rednic-este-noul-antrenor-echipei-viitorul.htm'

$ws.Range("A16").Value = 'Efectele pandemiei de COVID-19: Unul din 33 de locuitori ai lumii, număr record, va avea nevoie de ajutor pentru a supraviețui'
$ws.Range("B16").Value = 'https://www.hotnews.ro/stiri-international-24452405-efectele-pandemiei-covid-19-unul-din-33-locuitori-lumii-numar-record-avea-nevoie-ajutor-pentru-supravietui.htm'

$ws.Range("A17").Value = 'În Europa de Est publicul e menținut în afara consultărilor privind planul de redresare verde / În România, discuțiile s-au purtat în spatele ușilor închise (Emerging Europe)'
$ws.Range("B17").Value = 'https://www.hotnews.ro/stiri-international-24453618-europa-est-publicul-mentinut-afara-consultarilor-privind-planul-redresare-verde-romania-discutiile-purtat-spatele-usilor-inchise-emerging-europe.htm'

$ws.Range("A18").Value = 'Reacție dură a Patriarhiei după ce CE a propus desfășurarea online a slujbelor: Nimeni onest sau ancorat în realitatea socială profundă nu poate dori sau impune acest lucru'
$ws.Range("B18").Value = 'https://www.hotnews.ro/stiri-esential-24454472-reactie-dura-patriarhiei-dupa-propus-desfasurarea-online-slujbelor-nimeni-onest-sau-ancorat-realitatea-sociala-profunda-nu-poate-dori-sau-impune-acest-lucru.htm'

$ws.Range("A19").Value = 'Klaus Iohannis a decorat mai mulți doctori și asistente medicale, de Ziua Națională a României: Dvs, cei din linia întâi, luptați neobosit pentru viață, de multe ori cu riscul propriei siguranțe VIDEO'
$ws.Range("B19").Value = 'https://www.hotnews.ro/stiri-esential-24453142-klaus-iohannis-decorat-mai-multi-doctori-asistente-medicale-ziua-nationala-romaniei-video.htm'

$ws.Range("A20").Value = 'Giuliani ar fi discutat cu Trump despre posibilă grațiere - presă'
$ws.Range("B20").Value = 'https://www.hotnews.ro/stiri-international-24453261-giuliani-spus-discutat-posibila-gratiere-trump.htm'

$ws.Range("A21").Value = 'LISTĂ: Antreprenori IT români care s-au remarcat în anul 2020'
$ws.Range("B21").Value = 'https://economie.hotnews.ro/stiri-it-24453329-lista-antreprenori-romani-care-remarcat-anul-2020.htm'
